$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.586.42'
$ws.Range('E2').Value = '  -5.89%  '
$ws.Range('D3').Value = '3.294.34'
$ws.Range('E3').Value = '  -7.03%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '555.89'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -5.41%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '180.45'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -7.69%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.587'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.96%  '
$ws.Range('D9').Value = '3.291.81'
$ws.Range('E9').Value = '  -6.71%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.188'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -9.02%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.586'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -6.11%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '47.42'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -10.08%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000264'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -8.88%  '
$ws.Range('B14').Value = 'BitcoinCash'
$ws.Range('C14').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '629.54'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.44%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.54'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -7.89%  '
$ws.Range('B16').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C16').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D16').Value = '3.812.35'
$ws.Range('E16').Value = '  -7.18%  '
$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '17.95'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.67%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '65.480.38'
$ws.Range('E18').Value = '  -6.04%  '
$ws.Range('E19').Value = '  -4.20%  '
$ws.Range('D20').Value = '3.281.47'
$ws.Range('E20').Value = '  -7.69%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.40'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -9.27%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.903'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -6.16%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '17.73'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.45%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '104.25'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.01'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -8.57%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.99'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -9.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.80%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.68'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -8.81%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.43'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -7.84%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.76'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -8.77%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '30.62'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -7.70%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.91'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -6.51%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.26'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -7.95%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.06'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -6.01%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.105'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.61%  '
$ws.Range('B36').Value = 'Bittensor'
$ws.Range('C36').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '539.04'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.08%  '
$ws.Range('D37').Value = '3.732.43'
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('E38').Value = '  +0.20%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '56.78'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -7.97%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.45'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.65%  '
$ws.Range('B41').Value = 'CoreDAO'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.50'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +29.61%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.73'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -8.10%  '
$ws.Range('B43').Value = 'PEPE'
$ws.Range('C43').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D43').Value = '0.0₃0709'
$ws.Range('E43').Value = '  -12.53%  '
$ws.Range('B44').Value = 'Kaspa'
$ws.Range('C44').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.126'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -6.41%  '
$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.339'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -8.98%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '31.85'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -8.97%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0412'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -8.44%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.21'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -6.09%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.129'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -5.33%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.59'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -9.64%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.998'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.33%  '
